$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix 116: correct homophone typo 鑒賞 -> 鑑賞 ---
$ws.Range("L116").Value = '新城治安局因故需派出警員裝扮成貴族小姐進行潛入搜查，需要精通服裝搭配和香氛鑑賞的禁閉者支援。'

# --- Insert 4 new dispatch rows before the old row 174 (shifts it down to row 178) ---
$ws.Rows("174:177").Insert()

# Row 174: ID 1054001 - Friends in High Places
$ws.Range("A174").Value = 1054001
$ws.Range("B174").Value = 'Purple'
$ws.Range("C174").Value = 'Friends in High Places'
$ws.Range("D174").Value = '社会で頼れるのは友'
$ws.Range("E174").Value = '든든한 친구'
$ws.Range("F174").Value = '出门靠朋友'
$ws.Range("G174").Value = '出門靠朋友'
$ws.Range("H174").Value = 'According to intelligence reports, an Eastside tycoon is covertly financing a Mania weapons smuggling ring. The Bureau will dispatch undercover agents to infiltrate their network for further investigation. To gain access, an introduction from a high-society insider will be required.'
$ws.Range("I174").Value = '情報筋の話によると、ニューシティのとある富豪が狂瞳武器の密輸組織に密かに資金提供をしているという。管理局はその社交圏に作戦メンバーを潜入させる予定で、現在その富豪を紹介してくれる上流階級の名士を探している。'
$ws.Range("J174").Value = '정보원에 따르면 신성의 어느 한 부자가 암암리에 변이 무기 밀수 조직을 지원한다고 한다. 이에 따라 관리국은 작전 요원을 해당 인물의 사교계에 잠입시켜 조사를 진행할 예정이며, 현재 상류층 내에서 명망 있는 인물의 소개가 필요한 상황이다.'
$ws.Range("K174").Value = '据线报称，新城某富豪暗中资助狂厄武器走私组织，管理局将派出行动人员潜伏进入其社交圈进行调查，现需上流阶层有名望的人士为其引荐。'
$ws.Range("L174").Value = '據線報稱，新城某富豪暗中資助狂厄武器走私組織，管理局將派出行動人員潛伏進入其社交圈進行調查，現需上流階層有名望的人士為其引薦。'
$ws.Range("M174").Value = 'Parfait'
$ws.Range("N174").Value = 'パフェ'
$ws.Range("O174").Value = '파르·페'
$ws.Range("P174").Value = '芭·菲'
$ws.Range("Q174").Value = '芭·菲'
$ws.Range("R174").Value = ''
$ws.Range("S174").Value = ''
$ws.Range("T174").Value = ''
$ws.Range("U174").Value = ''
$ws.Range("V174").Value = ''
$ws.Range("W174").Value = 'Infected Cyst'
$ws.Range("X174").Value = '感染された嚢胞'
$ws.Range("Y174").Value = '감염된 낭포'
$ws.Range("Z174").Value = '感染囊胞'
$ws.Range("AA174").Value = '感染囊胞'
$ws.Range("AB174").Value = "'" + '1.0'
$ws.Range("AC174").Value = 'Organic Cyst'
$ws.Range("AD174").Value = '原生嚢胞'
$ws.Range("AE174").Value = '원시적 낭포'
$ws.Range("AF174").Value = '原生囊胞'
$ws.Range("AG174").Value = '原生囊胞'
$ws.Range("AH174").Value = "'" + '1.0'

# Row 175: ID 1054002 - Deadly Dance Party
$ws.Range("A175").Value = 1054002
$ws.Range("B175").Value = 'Purple'
$ws.Range("C175").Value = 'Deadly Dance Party'
$ws.Range("D175").Value = '死を招くダンス大会'
$ws.Range("E175").Value = '죽음의 무도회'
$ws.Range("F175").Value = '夺命舞会'
$ws.Range("G175").Value = '奪命舞會'
$ws.Range("H175").Value = 'Eastside is currently hosting a dance competition. Since the event began, M-values within the venue have been fluctuating abnormally. FAC has requested that the Bureau deploy a Sinner skilled in dance to accompany operatives for undercover investigation.'
$ws.Range("I175").Value = 'ニューシティでは最近ダンス大会が開催されているが、大会開始以来、会場内のM値が異常な波動を見せ続けている。FACは管理局に対し、ダンスに長けたコンビクトを作戦メンバーと共に潜入調査へ派遣するよう要請した。'
$ws.Range("J175").Value = '최근 신성에서 무용 대회를 개최하고 있다. 대회 시작 이후 경기장 내 M 수치가 지속적으로 비정상적인 파동을 보이고 있다. 이에 FAC는 관리국에 무용에 능한 수감자를 파견해 작전 요원과 함께 잠입 조사를 해달라고 요청했다.'
$ws.Range("K175").Value = '新城近日正在举办舞蹈大会，自开赛以来，场馆内M值一直波动异常，FAC向管理局申请派出一位擅长舞蹈的禁闭者随行动人员潜入调查。'
$ws.Range("L175").Value = '新城近日正在舉辦舞蹈大會，自開賽以來，場館內M值一直波動異常，FAC向管理局申請派出一位擅長舞蹈的禁閉者隨行動人員潛入調查。'
$ws.Range("M175").Value = 'Jasmine'
$ws.Range("N175").Value = 'ジャスミン'
$ws.Range("O175").Value = '자스민'
$ws.Range("P175").Value = '茉莉'
$ws.Range("Q175").Value = '茉莉'
$ws.Range("R175").Value = ''
$ws.Range("S175").Value = ''
$ws.Range("T175").Value = ''
$ws.Range("U175").Value = ''
$ws.Range("V175").Value = ''
$ws.Range("W175").Value = 'Ice Crystal'
$ws.Range("X175").Value = '氷晶'
$ws.Range("Y175").Value = '얼음 결정'
$ws.Range("Z175").Value = '冰晶'
$ws.Range("AA175").Value = '冰晶'
$ws.Range("AB175").Value = "'" + '1.0'
$ws.Range("AC175").Value = 'Ice Stone Concentrate'
$ws.Range("AD175").Value = '氷の精鉱'
$ws.Range("AE175").Value = '정교한 얼음석 광석'
$ws.Range("AF175").Value = '冰石精矿'
$ws.Range("AG175").Value = '冰石精礦'
$ws.Range("AH175").Value = "'" + '1.0'

# Row 176: ID 1054003 - Hidden Attraction
$ws.Range("A176").Value = 1054003
$ws.Range("B176").Value = 'Green'
$ws.Range("C176").Value = 'Hidden Attraction'
$ws.Range("D176").Value = '不人気スポット'
$ws.Range("E176").Value = '비인기 여행지'
$ws.Range("F176").Value = '冷门景点'
$ws.Range("G176").Value = '冷門景點'
$ws.Range("H176").Value = 'To raise residents'' awareness of Eastside''s public security institutions, the Public Security Bureau recently organized a one-day tour event. However, participation has been low, and a popular influencer is now needed to boost engagement.'
$ws.Range("I176").Value = 'ニューシティの治安管理機関に対する市民の理解を深めるため、治安局は最近一日見学イベントの開催を発表した。しかし応募者が少ないため、現在イベントの宣伝を担当する人気インフルエンサーを求めている。'
$ws.Range("J176").Value = '신성 치안 관리 기관에 대한 주민들의 이해를 높이기 위해 치안국에서 최근 일일 견학 프로그램을 기획했으나, 신청 인원이 매우 저조한 상황이다. 이에 따라 본 프로그램의 홍보를 맡아줄 인기 인플루언서의 도움이 필요하다.'
$ws.Range("K176").Value = '为提升居民对新城治安管理机构的了解，治安局近期开展了一日参观活动，但报名者寥寥，现需一名人气博主为活动进行宣传。'
$ws.Range("L176").Value = '為提升居民對新城治安管理機構的瞭解，治安局近期開展了一日參觀活動，但報名者寥寥，現需一名人氣網紅為活動進行宣傳。'
$ws.Range("M176").Value = 'Thalia'
$ws.Range("N176").Value = 'タリア'
$ws.Range("O176").Value = '탈리아'
$ws.Range("P176").Value = '塔利娅'
$ws.Range("Q176").Value = '塔利婭'
$ws.Range("R176").Value = ''
$ws.Range("S176").Value = ''
$ws.Range("T176").Value = ''
$ws.Range("U176").Value = ''
$ws.Range("V176").Value = ''
$ws.Range("W176").Value = 'Ice Stone Concentrate'
$ws.Range("X176").Value = '氷の精鉱'
$ws.Range("Y176").Value = '정교한 얼음석 광석'
$ws.Range("Z176").Value = '冰石精矿'
$ws.Range("AA176").Value = '冰石精礦'
$ws.Range("AB176").Value = "'" + '1.0'
$ws.Range("AC176").Value = 'Ice Stone Raw Ore'
$ws.Range("AD176").Value = '氷の原鉱'
$ws.Range("AE176").Value = '거친 얼음석 광석'
$ws.Range("AF176").Value = '冰石粗矿'
$ws.Range("AG176").Value = '冰石粗礦'
$ws.Range("AH176").Value = "'" + '1.0'

# Row 177: ID 1054004 - Top Influencer
$ws.Range("A177").Value = 1054004
$ws.Range("B177").Value = 'Blue'
$ws.Range("C177").Value = 'Top Influencer'
$ws.Range("D177").Value = 'トップ100インフルエンサー'
$ws.Range("E177").Value = '인플루언서 TOP 100'
$ws.Range("F177").Value = '百大博主'
$ws.Range("G177").Value = '百大網紅'
$ws.Range("H177").Value = 'Recently, several online influencers in Eastside have mysteriously disappeared after attending industry gatherings. The Bureau urgently needs to deploy a Sinner with established recognition in the influencer community for undercover investigation.'
$ws.Range("I177").Value = '最近、ニューシティでは業界の集まりに招待された複数のインフルエンサーが謎の失踪を遂げる事件が相次いでいる。管理局は、インフルエンサー界隈である程度の知名度を持つコンビクトを潜入調査に派遣する必要がある。'
$ws.Range("J177").Value = '최근, 신성의 여러 유명 블로거가 업계 모임에 초청받은 후 실종되었다. 관리국에서 신속히 인플루언서 업계에서 유명한 수감자를 파견해 잠복 조사를 진행해야 한다.'
$ws.Range("K177").Value = '近期，新城多名网络博主在受邀参与业内聚会后神秘失踪，急需管理局派出在网红界小有名气的禁闭者卧底调查。'
$ws.Range("L177").Value = '近期，新城多名網紅在受邀參與業內聚會後神秘失蹤，急需管理局派出在網紅界小有名氣的禁閉者臥底調查。'
$ws.Range("M177").Value = 'Thalia'
$ws.Range("N177").Value = 'タリア'
$ws.Range("O177").Value = '탈리아'
$ws.Range("P177").Value = '塔利娅'
$ws.Range("Q177").Value = '塔利婭'
$ws.Range("R177").Value = 'Nino'
$ws.Range("S177").Value = 'ニノ'
$ws.Range("T177").Value = '니노'
$ws.Range("U177").Value = '妮诺'
$ws.Range("V177").Value = '妮諾'
$ws.Range("W177").Value = 'Ice Stone Concentrate'
$ws.Range("X177").Value = '氷の精鉱'
$ws.Range("Y177").Value = '정교한 얼음석 광석'
$ws.Range("Z177").Value = '冰石精矿'
$ws.Range("AA177").Value = '冰石精礦'
$ws.Range("AB177").Value = "'" + '1.0'
$ws.Range("AC177").Value = 'Ice Stone Concentrate'
$ws.Range("AD177").Value = '氷の精鉱'
$ws.Range("AE177").Value = '정교한 얼음석 광석'
$ws.Range("AF177").Value = '冰石精矿'
$ws.Range("AG177").Value = '冰石精礦'
$ws.Range("AH177").Value = "'" + '1.0'

# Row 178: ID 1099001 - Dream Interpretation
$ws.Range("A178").Value = 1099001
$ws.Range("B178").Value = 'Green'
$ws.Range("C178").Value = 'Dream Interpretation'
$ws.Range("D178").Value = '夢の解析'
$ws.Range("E178").Value = '꿈의 해석'
$ws.Range("F178").Value = '梦的解析'
$ws.Range("G178").Value = '夢的解析'
$ws.Range("H178").Value = 'Residents in some areas of Eastside often suffer from nightmares. It is necessary to find out why.'
$ws.Range("I178").Value = 'ニューシティの一部エリアの住民がよくナイトメアにうなされている。具体的な原因を調査しなければならない。'
$ws.Range("J178").Value = '신성 일부 지역 주민들은 악몽을 자주꾼다. 구체적인 원인을 철저히 조사해야 한다.'
$ws.Range("K178").Value = '新城一些区域的居民经常做噩梦，需要查清具体原因。'
$ws.Range("L178").Value = '新城一些區域的居民經常做惡夢，需要查明具體原因。'
$ws.Range("M178").Value = 'Hecate'
$ws.Range("N178").Value = 'ヘカテー'
$ws.Range("O178").Value = '헤카테'
$ws.Range("P178").Value = '赫卡蒂'
$ws.Range("Q178").Value = '赫卡蒂'
$ws.Range("R178").Value = ''
$ws.Range("S178").Value = ''
$ws.Range("T178").Value = ''
$ws.Range("U178").Value = ''
$ws.Range("V178").Value = ''
$ws.Range("W178").Value = 'Arsenopyrite Concentrate'
$ws.Range("X178").Value = '毒砂の精鉱'
$ws.Range("Y178").Value = '정교한 독모래 광석'
$ws.Range("Z178").Value = '毒砂精矿'
$ws.Range("AA178").Value = '毒砂精礦'
$ws.Range("AB178").Value = "'" + '1.0'
$ws.Range("AC178").Value = 'Arsenopyrite Raw Ore'
$ws.Range("AD178").Value = '毒砂の原鉱'
$ws.Range("AE178").Value = '거친 독모래 광석'
$ws.Range("AF178").Value = '毒砂粗矿'
$ws.Range("AG178").Value = '毒砂粗礦'
$ws.Range("AH178").Value = "'" + '1.0'
